$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.238.99"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "3.891.25"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'483.74"
$ws.Range("E5").Value = "  -0.25%  "
$ws.Range("D6").Value = "'145.66"
$ws.Range("E6").Value = "  -1.64%  "
$ws.Range("E7").Value = "  +0.42%  "
$ws.Range("D8").Value = "'0.997"
$ws.Range("D9").Value = "'0.743"
$ws.Range("E9").Value = "  +2.79%  "
$ws.Range("D10").Value = "'0.182"
$ws.Range("E10").Value = "  +8.33%  "
$ws.Range("D11").Value = "'0.0000355"
$ws.Range("E11").Value = "  +0.50%  "
$ws.Range("D12").Value = "'43.16"
$ws.Range("E12").Value = "  +1.65%  "
$ws.Range("D13").Value = "'10.54"
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("D14").Value = "4.504.24"
$ws.Range("E14").Value = "  -1.12%  "
$ws.Range("D15").Value = "3.909.33"
$ws.Range("E15").Value = "  -1.05%  "
$ws.Range("E16").Value = "  -2.82%  "
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").Value = "'19.99"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("D20").Value = "68.219.71"
$ws.Range("E20").Value = "  -0.44%  "
$ws.Range("D21").Value = "'430.32"
$ws.Range("E21").Value = "  -0.18%  "
$ws.Range("D22").Value = "'3.59"
$ws.Range("E22").Value = "  +7.46%  "
$ws.Range("E23").Value = "  +2.20%  "
$ws.Range("E24").Value = "  +2.86%  "
$ws.Range("D25").Value = "'12.26"
$ws.Range("E25").Value = "  +16.89%  "
$ws.Range("E26").Value = "  +2.74%  "
$ws.Range("D27").Value = "'11.04"
$ws.Range("E27").Value = "  -2.17%  "
$ws.Range("D28").Value = "'37.36"
$ws.Range("E28").Value = "  -1.94%  "
$ws.Range("D29").Value = "'5.67"
$ws.Range("E29").Value = "  -3.58%  "
$ws.Range("D30").Value = "'711.90"
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("D31").Value = "'13.48"
$ws.Range("E31").Value = "  +1.90%  "
$ws.Range("D32").Value = "'0.130"
$ws.Range("E32").Value = "  +0.55%  "
$ws.Range("D33").Value = "'2.92"
$ws.Range("E33").Value = "  +3.03%  "
$ws.Range("D34").Value = "0.0₃0883"
$ws.Range("E34").Value = "  -1.64%  "
$ws.Range("D35").Value = "'61.80"
$ws.Range("E35").Value = "  +5.86%  "
$ws.Range("D36").Value = "'6.06"
$ws.Range("E36").Value = "  +10.37%  "
$ws.Range("D37").Value = "'40.93"
$ws.Range("E37").Value = "  -1.48%  "
$ws.Range("E38").Value = "  -3.29%  "
$ws.Range("B39").Value = "TheGraph"
$ws.Range("C39").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D39").Value = "'0.396"
$ws.Range("E39").Value = "  +15.65%  "
$ws.Range("E40").Value = "  +6.66%  "
$ws.Range("B41").Value = "Dai"
$ws.Range("C41").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D41").Value = "'0.998"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").Value = "'3.01"
$ws.Range("E42").Value = "  +5.44%  "
$ws.Range("E43").Value = "  +2.62%  "
$ws.Range("E44").Value = "  -3.45%  "
$ws.Range("E45").Value = "  +1.60%  "
$ws.Range("E46").Value = "  +3.35%  "
$ws.Range("E47").Value = "  -0.08%  "
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("E49").Value = "  -2.24%  "
$ws.Range("D50").Value = "'144.75"
$ws.Range("E50").Value = "  -1.51%  "
$ws.Range("D51").Value = "0.0₆0337"
$ws.Range("E51").Value = "  +24.47%  "
